$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all updated cells retain their original text formatting (no auto numeric conversion)
$textCells = @("D2","D3","D5","D7","D8","D9","D10","D12","D13","D14","D15","D16","D17","D19","D20","D21","D22","D23","D25","D26","D27","D28","D29","D30","D33","D35","D36","D37","D38","D39","D40","D41","D42","D44","D45","D46","D47","D48","D50","D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$updates = @{
    "D2" = '30.650.08'
    "E2" = '  +0.71%  '
    "D3" = '2.125.74'
    "E3" = '  +0.90%  '
    "E4" = '  +0.79%  '
    "D5" = '352.61'
    "E5" = '  +5.43%  '
    "E6" = '  +0.65%  '
    "D7" = '0.5281'
    "E7" = '  +0.97%  '
    "D8" = '0.4546'
    "E8" = '  -0.19%  '
    "D9" = '54.14'
    "E9" = '  +1.64%  '
    "D10" = '0.09095'
    "E10" = '  +1.73%  '
    "E11" = '  +0.53%  '
    "D12" = '24.67'
    "E12" = '  +1.64%  '
    "D13" = '2.136.00'
    "E13" = '  +1.58%  '
    "D14" = '6.870'
    "E14" = '  +0.20%  '
    "D15" = '8.127'
    "E15" = '  +0.93%  '
    "D16" = '102.59'
    "E16" = '  +6.22%  '
    "D17" = '0.00001177'
    "E17" = '  +3.05%  '
    "E18" = '  +0.64%  '
    "D19" = '0.06714'
    "E19" = '  +0.85%  '
    "D20" = '19.49'
    "E20" = '  +1.29%  '
    "D21" = '1.010'
    "E21" = '  +0.60%  '
    "D22" = '6.358'
    "E22" = '  +0.30%  '
    "D23" = '30.732.57'
    "E23" = '  +0.71%  '
    "E24" = '  +3.19%  '
    "D25" = '2.389'
    "E25" = '  +1.14%  '
    "D26" = '2.386.14'
    "E26" = '  +1.46%  '
    "D27" = '22.54'
    "E27" = '  +1.12%  '
    "D28" = '2.569'
    "E28" = '  +1.02%  '
    "D29" = '164.82'
    "E29" = '  +1.26%  '
    "D30" = '136.52'
    "E30" = '  +2.51%  '
    "E32" = '  +0.90%  '
    "D33" = '1.670'
    "E33" = '  +0.07%  '
    "E34" = '  +0.13%  '
    "D35" = '4.020'
    "E35" = '  +1.93%  '
    "D36" = '6.207'
    "E36" = '  +8.62%  '
    "D37" = '10.31'
    "E37" = '  -1.94%  '
    "D38" = '0.02656'
    "E38" = '  +2.56%  '
    "D39" = '0.06895'
    "E39" = '  +0.58%  '
    "D40" = '0.2327'
    "E40" = '  +1.04%  '
    "D41" = '12.56'
    "E41" = '  -1.21%  '
    "D42" = '0.6924'
    "E42" = '  +0.47%  '
    "E43" = '  +2.16%  '
    "D44" = '14.85'
    "E44" = '  +5.86%  '
    "D45" = '2.342'
    "E45" = '  +0.90%  '
    "D46" = '0.6460'
    "E46" = '  +1.19%  '
    "D47" = '3.760'
    "E47" = '  +2.63%  '
    "D48" = '0.00000000366'
    "E48" = '  +5.51%  '
    "E49" = '  +0.41%  '
    "B50" = 'WOONetwork'
    "C50" = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
    "D50" = '0.3386'
    "E50" = '  -0.73%  '
    "B51" = 'Aave'
    "C51" = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    "D51" = '83.06'
    "E51" = '  -0.59%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
